$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values per row/column, as produced by the refreshed NATMI script run.
$updates = @{
    2 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 136.544502; "N" = 409.633506; "O" = 0.9681180443787725; "P" = 0.9681180443787725; "Q" = 359.075537363132; "R" = 3231.679836268188; "S" = 0.0692856700765254; "T" = 0.0692856700765254 }
    3 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "O" = 0.002758738216274633; "P" = 0.002758738216274633; "Q" = 1.023217585091778; "R" = 9.208958265826; "S" = 0.0001974356608578232; "T" = 0.0001974356608578232 }
    4 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 4.025396333333334; "N" = 12.076189; "O" = 0.0285405766544606; "P" = 0.02854057665446059; "Q" = 10.58571623404689; "R" = 95.27144610642202; "S" = 0.002042574239119408; "T" = 0.002042574239119408 }
    5 = @{ "G" = 2.629732666666667; "H" = 7.889198; "I" = 0.07156737804735891; "J" = 0.07156737804735891; "M" = 0.08217633333333334; "N" = 0.246529; "O" = 0.0005826407504923545; "P" = 0.0005826407504923544; "Q" = 0.2161017881935556; "R" = 1.944916093742; "S" = 0.00004169807085628325; "T" = 0.00004169807085628324 }
    6 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 136.544502; "N" = 409.633506; "O" = 0.9681180443787725; "P" = 0.9681180443787725; "Q" = 2475.098174909522; "R" = 22275.8835741857; "S" = 0.4775842899605972; "T" = 0.4775842899605972 }
    7 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "O" = 0.002758738216274633; "P" = 0.002758738216274633; "S" = 0.00136091878449815; "T" = 0.00136091878449815 }
    8 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 4.025396333333334; "N" = 12.076189; "O" = 0.0285405766544606; "P" = 0.02854057665446059; "Q" = 72.96706181491524; "R" = 656.7035563342371; "S" = 0.01407941016669417; "T" = 0.01407941016669416 }
    9 = @{ "I" = 0.493312042610523; "J" = 0.493312042610523; "M" = 0.08217633333333334; "N" = 0.246529; "O" = 0.0005826407504923545; "P" = 0.0005826407504923544; "Q" = 1.489583906161889; "R" = 13.406255155457; "S" = 0.0002874236987335115; "T" = 0.0002874236987335115 }
    10 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 136.544502; "N" = 409.633506; "O" = 0.9681180443787725; "P" = 0.9681180443787725; "Q" = 1031.000608808146; "R" = 9279.005479273315; "S" = 0.198937439612706; "T" = 0.198937439612706 }
    11 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "O" = 0.002758738216274633; "P" = 0.002758738216274633; "Q" = 2.937927659789222; "R" = 26.441348938103; "S" = 0.0005668898751490233; "T" = 0.0005668898751490235 }
    12 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 4.025396333333334; "N" = 12.076189; "O" = 0.0285405766544606; "P" = 0.02854057665446059; "Q" = 30.39438431846012; "R" = 273.549458866141; "S" = 0.005864769567797817; "T" = 0.005864769567797817 }
    13 = @{ "G" = 7.550656333333333; "H" = 22.651969; "I" = 0.2054888252189962; "J" = 0.2054888252189962; "M" = 0.08217633333333334; "N" = 0.246529; "O" = 0.0005826407504923545; "P" = 0.0005826407504923544; "Q" = 0.6204852517334445; "R" = 5.584367265601; "S" = 0.0001197261633433882; "T" = 0.0001197261633433882 }
    14 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 136.544502; "N" = 409.633506; "O" = 0.9681180443787725; "P" = 0.9681180443787725; "Q" = 1152.133105293236; "R" = 10369.19794763912; "S" = 0.2223106447289439; "T" = 0.2223106447289439 }
    15 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "O" = 0.002758738216274633; "P" = 0.002758738216274633; "Q" = 3.283105450066445; "R" = 29.547949050598; "S" = 0.0006334938957696363; "T" = 0.0006334938957696363 }
    16 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 4.025396333333334; "N" = 12.076189; "O" = 0.0285405766544606; "P" = 0.02854057665446059; "Q" = 33.96542745865623; "R" = 305.6888471279061; "S" = 0.006553822680849209; "T" = 0.006553822680849208 }
    17 = @{ "G" = 8.437784666666667; "H" = 25.313354; "I" = 0.2296317541231219; "J" = 0.2296317541231219; "M" = 0.08217633333333334; "N" = 0.246529; "O" = 0.0005826407504923545; "P" = 0.0005826407504923544; "Q" = 0.693386205362889; "R" = 6.240475848266; "S" = 0.0001337928175591715; "T" = 0.0001337928175591715 }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
